$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 11:10"

# Row 6 - India: updated case numbers
$ws.Range("B6").Value = 1808128
$ws.Range("C6").Value = 3426
$ws.Range("D6").Value = 1188650
$ws.Range("E6").Value = 581277
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 38201

# Row 26 - Indonesia: updated case numbers
$ws.Range("B26").Value = 113134
$ws.Range("C26").Value = 1679
$ws.Range("D26").Value = 70237
$ws.Range("E26").Value = 37595
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 66
$ws.Range("H26").Value = 5302

# Row 28 - Filipinas: updated case numbers
$ws.Range("B28").Value = 106330
$ws.Range("C28").Value = 3226
$ws.Range("D28").Value = 65821
$ws.Range("E28").Value = 38405
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 45
$ws.Range("H28").Value = 2104

# Rows 36/37 - Israel overtakes Ucrania in the sorted ranking.
# Row 36 becomes Israel (with Israel's new figures), row 37 becomes Ucrania
# (with Ucrania's figures, which stay numerically the same).
$ws.Range("A36").Value = "Israel"
$ws.Range("B36").Value = 73231
$ws.Range("C36").Value = 416
$ws.Range("D36").Value = 47523
$ws.Range("E36").Value = 25167
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 541

$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 73158
$ws.Range("C37").Value = 990
$ws.Range("D37").Value = 39876
$ws.Range("E37").Value = 31544
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 13
$ws.Range("H37").Value = 1738

# Row 49 - Polonia: updated case numbers
$ws.Range("B49").Value = 47469
$ws.Range("C49").Value = 575
$ws.Range("D49").Value = 34881
$ws.Range("E49").Value = 10856
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 1732

# Row 93 - Finlandia: updated case numbers
$ws.Range("B93").Value = 7466
$ws.Range("C93").Value = 13
$ws.Range("D93").Value = 6950
$ws.Range("E93").Value = 187

# Row 112 - Hong Kong: updated case numbers
$ws.Range("B112").Value = 3592
$ws.Range("C112").Value = 80
$ws.Range("D112").Value = 2037
$ws.Range("E112").Value = 1518

# Row 126 - Eslovenia: updated case numbers
$ws.Range("B126").Value = 2181
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 1826
$ws.Range("E126").Value = 233
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 122
